# Table 5.6 special processing
# Adds a new worksheet "5.6 Annual summaries" at the end of the workbook,
# containing a flattened (row/group/item/label/unit) summary table used
# for downstream processing of table 5.6.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last sheet ("Foglio8"), becoming the new
# last tab and the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "5.6 Annual summaries"

# Column widths (matches sheet8/sheet9 style columns). The host's column
# width is stored at pixel granularity, so the inputs below are chosen to
# land as close as possible to the source file's exact character widths
# (22.21875 / 22) after that internal rounding.
$ws.Columns.Item(2).ColumnWidth = 21.25
$ws.Columns.Item(3).ColumnWidth = 21.1

# Zoom as in the source file
$ws.Activate()
$excel.ActiveWindow.Zoom = 190

# Header row
$ws.Range("A1").Value = "row"
$ws.Range("B1").Value = "group"
$ws.Range("C1").Value = "item"
$ws.Range("D1").Value = "label"
$ws.Range("E1").Value = "unit"

# Data rows: group / item pairs, one per generator-category x measure.
$ws.Range("B2").Value = "Major power producers"
$ws.Range("C2").Value = "Fuel used "
$ws.Range("E2").Value = "GWh"

$ws.Range("B3").Value = "Major power producers"
$ws.Range("C3").Value = "Generation"
$ws.Range("E3").Value = "GWh"

$ws.Range("B4").Value = "Major power producers"
$ws.Range("C4").Value = "Used on works"
$ws.Range("E4").Value = "GWh"

$ws.Range("B5").Value = "Major power producers"
$ws.Range("C5").Value = "Supplied (gross)"
$ws.Range("E5").Value = "GWh"

$ws.Range("B6").Value = "Major power producers"
$ws.Range("C6").Value = "Used in pumping"
$ws.Range("E6").Value = "GWh"

$ws.Range("B7").Value = "Major power producers"
$ws.Range("C7").Value = "Supplied (net)"
$ws.Range("E7").Value = "GWh"

$ws.Range("B8").Value = "Other generators"
$ws.Range("C8").Value = "Fuel used "
$ws.Range("E8").Value = "GWh"

$ws.Range("B9").Value = "Other generators"
$ws.Range("C9").Value = "Generation "
$ws.Range("E9").Value = "GWh"

$ws.Range("B10").Value = "Other generators"
$ws.Range("C10").Value = "Used on works"
$ws.Range("E10").Value = "GWh"

$ws.Range("B11").Value = "Other generators"
$ws.Range("C11").Value = "Supplied "
$ws.Range("E11").Value = "GWh"

$ws.Range("B12").Value = "All generating companies"
$ws.Range("C12").Value = "Fuel used"
$ws.Range("E12").Value = "GWh"

$ws.Range("B13").Value = "All generating companies"
$ws.Range("C13").Value = "Generation "
$ws.Range("E13").Value = "GWh"

$ws.Range("B14").Value = "All generating companies"
$ws.Range("C14").Value = "Used on works"
$ws.Range("E14").Value = "GWh"

$ws.Range("B15").Value = "All generating companies"
$ws.Range("C15").Value = "Supplied (gross)"
$ws.Range("E15").Value = "GWh"

$ws.Range("B16").Value = "All generating companies"
$ws.Range("C16").Value = "Used in pumping"
$ws.Range("E16").Value = "GWh"

$ws.Range("B17").Value = "All generating companies"
$ws.Range("C17").Value = "Supplied (net)"
$ws.Range("E17").Value = "GWh"

# Column A: running row counter 0..15, built with a relative fill-down formula
# (row 3 gets its own formula, rows 4:17 share it) to match Excel's native
# shared-formula behaviour when filling a column down.
$ws.Range("A2").Value = 0
$ws.Range("A3").Formula = "=1+A2"
$ws.Range("A4:A17").Formula = "=1+A3"

# Column D: concatenated "group item" label, same fill-down shared-formula
# pattern as column A.
$ws.Range("D2").Formula = "=_xlfn.CONCAT(B2,`" `",C2)"
$ws.Range("D3:D17").Formula = "=_xlfn.CONCAT(B3,`" `",C3)"

# Selection / view state to match the source file.
$ws.Range("A18").Select() | Out-Null
